$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data (2021-09-09, serial 44448) was added to the top of
# the "Palta" (avocado) block. The existing rows 199-207 (weeks of 2020-12-02,
# 2021-02-11 and 2021-07-22) shift down to rows 205-213, and six brand-new
# rows are inserted at 199-204 for the new week.
$ws.Rows("199:204").Insert()

# New rows to populate (row, Fecha, Variedad, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion, Precio $/Kg)
$data = @(
    @{Row=199; D=44448; K="Edranol"; L="Especial";     M=200; N=2300; O=2400; P=2350; Q="`$/kilo (en caja de 17 kilos)"; S=2350},
    @{Row=200; D=44448; K="Edranol"; L="Primera";      M=240; N=2100; O=2200; P=2150; Q="`$/kilo (en caja de 17 kilos)"; S=2150},
    @{Row=201; D=44448; K="Edranol"; L="Segunda";      M=200; N=1900; O=2000; P=1950; Q="`$/kilo (en caja de 17 kilos)"; S=1950},
    @{Row=202; D=44448; K="Hass";    L="1a nueva(o)";  M=600; N=2500; O=2600; P=2550; Q="`$/kilo (en caja de 17 kilos)"; S=2550},
    @{Row=203; D=44448; K="Hass";    L="2a nueva(o)";  M=400; N=2200; O=2300; P=2250; Q="`$/kilo (en caja de 17 kilos)"; S=2250},
    @{Row=204; D=44448; K="Hass";    L="3a nueva (o)"; M=200; N=2000; O=2100; P=2050; Q="`$/kilo (en caja de 17 kilos)"; S=2050}
)

foreach ($d in $data) {
    $r = $d.Row
    $ws.Cells.Item($r, 1).Value  = 2
    $ws.Cells.Item($r, 2).Value  = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($r, 3).Value  = "Coquimbo"
    $ws.Cells.Item($r, 4).Value  = $d.D
    $ws.Cells.Item($r, 5).Value  = 4
    $ws.Cells.Item($r, 6).Value  = "Fruta"
    $ws.Cells.Item($r, 7).Value  = 100106
    $ws.Cells.Item($r, 8).Value  = "Oleaginosos"
    $ws.Cells.Item($r, 9).Value  = 100106002
    $ws.Cells.Item($r, 10).Value = "Palta"
    $ws.Cells.Item($r, 11).Value = $d.K
    $ws.Cells.Item($r, 12).Value = $d.L
    $ws.Cells.Item($r, 13).Value = $d.M
    $ws.Cells.Item($r, 14).Value = $d.N
    $ws.Cells.Item($r, 15).Value = $d.O
    $ws.Cells.Item($r, 16).Value = $d.P
    $ws.Cells.Item($r, 17).Value = $d.Q
    $ws.Cells.Item($r, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 19).Value = $d.S
    $ws.Cells.Item($r, 20).Value = 1
}
